# add flexible column entry to generic parser
# Inserts a new "Mark" column before the existing "Stock" column (H) on the
# PIT-tagging template, matching header styling used by the other columns,
# and moves the trailing comments so they stay attached to the header they
# were originally documenting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at H - this shifts H:R to I:S, carrying
#    along cell values/styles/widths already in use.
$ws.Range("H1").EntireColumn.Insert()

# 2. The cell comments are anchored to their original row/col, so they do
#    NOT follow their header cell when the column shifts right. Re-home
#    each one, working from the rightmost column back towards H so a
#    not-yet-processed comment is never overwritten by one moving into its
#    old slot.
$commentCols = @("R", "P", "O", "N", "M", "L", "K", "J", "I", "H")
foreach ($col in $commentCols) {
    $src = $ws.Range($col + "3")
    if (-not ($src.Comment -eq $null)) {
        $text = $src.Comment.Text()
        $src.Comment.Delete()
        $destCol = [char]([int][char]$col + 1)
        $dest = $ws.Range($destCol + "3")
        $dest.AddComment($text)
    }
}

# 3. Populate the new "Mark" header cell and give it the same look as the
#    other header cells (bold font + bottom border already carried over
#    from the column insert; add the slightly-darker fill + centering used
#    for this new column).
$h3 = $ws.Range("H3")
$h3.Value = "Mark"
$h3.Interior.Color = 10921638
$h3.HorizontalAlignment = -4108

# 4. Document the new column the same way the others are documented.
$h3.AddComment("Group mark, optional. `nEg. Adipose Clip`nMust match code in DB.")

# 5. Match the recorded selection left behind by the edit.
$ws.Range("H6").Select()
